# edit.ps1 - applies the "support for variable negation" commit to the
# Objeck Programming Language release-notes document.
#
# Summary of changes performed:
#   1. Remove the stray _GoBack bookmark that originally sat next to the
#      "10" in the date line.
#   2. Rewrite the "Finally, added support for multivariable declarations
#      and assignments [new]" bullet into "Added basic overdue language
#      features" (re-adding the _GoBack bookmark around "overdue ").
#   3. Add a new sub-bullet (list level 2) carrying the old multivariable-
#      declarations text: "Support for multivariable declarations and
#      assignments [new] ".
#   4. Add another new sub-bullet (list level 2): "Support for variable
#      negation [new]".
#   5. Merge the "Minor enhancements"/" and bug fixes"/" to the debugger
#      [new]"/" " runs of the last Highlights bullet into one run.

$d = $word.ActiveDocument

function Split-RunAt([int]$pos) {
    # Forces a run boundary at an absolute document character offset by
    # dropping a temporary bookmark there and immediately deleting it.
    $pt = $d.Range($pos, $pos)
    $d.Bookmarks.Add("__tmp_split__", $pt) | Out-Null
    $d.Bookmarks("__tmp_split__").Delete()
}

# ---------------------------------------------------------------------
# 1. Remove the old _GoBack bookmark (sitting next to "10" in the date).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Rewrite the "Finally, added support ..." bullet.
# ---------------------------------------------------------------------
$p8 = $d.Paragraphs(8)
$p8FullRange = $d.Range($p8.Range.Start, $p8.Range.End - 1)
$p8FullRange.Text = "Added basic overdue language features"

$p8Start = $d.Paragraphs(8).Range.Start
$bmStart = $p8Start + "Added basic ".Length
$bmEnd = $p8Start + "Added basic overdue ".Length
$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------
# 3. Insert the new "Support for multivariable declarations ..." bullet
#    right after paragraph 8, as a level-2 (ilvl=1) list item.
# ---------------------------------------------------------------------
$p8 = $d.Paragraphs(8)
$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs(9)
$p9.Range.ListFormat.ListLevelNumber = 2
$r9 = $p9.Range
$text9 = "Support for multivariable declarations and assignments [new] "
$r9.InsertBefore($text9)

$p9Start = $p9.Range.Start
$off1 = "S".Length
$off2 = $off1 + "upport for multivariable declarations and assignments ".Length
$off3 = $off2 + "[new]".Length
Split-RunAt ($p9Start + $off1)
Split-RunAt ($p9Start + $off2)
Split-RunAt ($p9Start + $off3)

# ---------------------------------------------------------------------
# 4. Insert the new "Support for variable negation [new]" bullet right
#    after the multivariable-declarations bullet, same list level.
# ---------------------------------------------------------------------
$p9 = $d.Paragraphs(9)
$p9.Range.InsertParagraphAfter()
$p10 = $d.Paragraphs(10)
$p10.Range.ListFormat.ListLevelNumber = 2
$r10 = $p10.Range
$text10 = "Support for variable negation [new]"
$r10.InsertBefore($text10)

$p10Start = $p10.Range.Start
$q1 = "S".Length
$q2 = $q1 + "upport for ".Length
$q3 = $q2 + "variable negation".Length
$q4 = $q3 + " ".Length
Split-RunAt ($p10Start + $q1)
Split-RunAt ($p10Start + $q2)
Split-RunAt ($p10Start + $q3)
Split-RunAt ($p10Start + $q4)

# ---------------------------------------------------------------------
# 5. Merge the "Minor enhancements ... [new] " bullet runs into one.
# ---------------------------------------------------------------------
$searchRange = $d.Range(0, $d.Content.End)
$searchRange.Find.Execute("Minor enhancements and bug fixes to the debugger [new] ", $true, $false, $false, $false, $false, $true, 1, $false, "Minor enhancements and bug fixes to the debugger [new] ", 2) | Out-Null
